$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.573.81"
$ws.Range("E2").Value = "  +6.02%  "

$ws.Range("D3").Value = "2.648.53"
$ws.Range("E3").Value = "  +9.89%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'512.62"
$ws.Range("E5").Value = "  +4.83%  "

$ws.Range("D6").Value = "'157.80"
$ws.Range("E6").Value = "  +2.38%  "

$ws.Range("D7").Value = "'0.993"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("D9").Value = "2.695.29"
$ws.Range("E9").Value = "  +10.94%  "

$ws.Range("E10").Value = "  +2.67%  "

$ws.Range("E11").Value = "  +5.20%  "

$ws.Range("E12").Value = "  +4.02%  "

$ws.Range("E13").Value = "  +1.06%  "

$ws.Range("D14").Value = "3.138.73"
$ws.Range("E14").Value = "  +10.10%  "

$ws.Range("D15").Value = "60.686.41"
$ws.Range("E15").Value = "  +6.26%  "

$ws.Range("E16").Value = "  +5.48%  "

$ws.Range("E17").Value = "  +5.29%  "

$ws.Range("D18").Value = "2.681.31"
$ws.Range("E18").Value = "  +10.35%  "

$ws.Range("E19").Value = "  +1.70%  "

$ws.Range("D20").Value = "'349.60"
$ws.Range("E20").Value = "  +7.97%  "

$ws.Range("D21").Value = "'10.56"
$ws.Range("E21").Value = "  +5.66%  "

$ws.Range("D22").Value = "'6.22"
$ws.Range("E22").Value = "  +4.10%  "

$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "'60.37"
$ws.Range("E24").Value = "  +3.71%  "

$ws.Range("E25").Value = "  +4.09%  "

$ws.Range("D26").Value = "2.794.83"
$ws.Range("E26").Value = "  +10.29%  "

$ws.Range("D27").Value = "'0.167"
$ws.Range("E27").Value = "  +3.80%  "

$ws.Range("D28").Value = "'0.995"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("D29").Value = "0.0₃0873"
$ws.Range("E29").Value = "  +11.55%  "

$ws.Range("D30").Value = "'7.56"
$ws.Range("E30").Value = "  +3.49%  "

$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  +6.01%  "

$ws.Range("D33").Value = "'157.38"
$ws.Range("E33").Value = "  +4.60%  "

$ws.Range("E34").Value = "  +3.40%  "

$ws.Range("D35").Value = "'5.75"
$ws.Range("E35").Value = "  +8.77%  "

$ws.Range("D36").Value = "'4.08"
$ws.Range("E36").Value = "  +9.46%  "

$ws.Range("E37").Value = "  +5.13%  "

$ws.Range("D38").Value = "'314.14"
$ws.Range("E38").Value = "  +17.01%  "

$ws.Range("E39").Value = "  +10.17%  "

$ws.Range("D40").Value = "'0.858"
$ws.Range("E40").Value = "  +1.25%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.78"
$ws.Range("E41").Value = "  +6.90%  "

$ws.Range("B42").Value = "SuiNetwork"
$ws.Range("C42").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D42").Value = "'0.838"
$ws.Range("E42").Value = "  +30.48%  "

$ws.Range("D43").Value = "'35.49"
$ws.Range("E43").Value = "  +3.91%  "

$ws.Range("D44").Value = "'0.649"
$ws.Range("E44").Value = "  +8.99%  "

$ws.Range("E45").Value = "  +8.51%  "

$ws.Range("E46").Value = "  -1.00%  "

$ws.Range("D47").Value = "'20.24"
$ws.Range("E47").Value = "  +16.16%  "

$ws.Range("D48").Value = "'0.991"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.91"
$ws.Range("E49").Value = "  +7.82%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.079.62"
$ws.Range("E50").Value = "  +10.86%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0237"
$ws.Range("E51").Value = "  +3.70%  "
